$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 23.699655
$ws.Range("H2").Value = 71.09896499999999
$ws.Range("I2").Value = 0.4841969272415696
$ws.Range("J2").Value = 0.4841969272415697
$ws.Range("M2").Value = 0.8366046666666667
$ws.Range("N2").Value = 2.509814
$ws.Range("O2").Value = 0.08025679986157715
$ws.Range("P2").Value = 0.08025679986157715
$ws.Range("Q2").Value = 19.82724197139
$ws.Range("R2").Value = 178.44517774251
$ws.Range("S2").Value = 0.03886009588321728
$ws.Range("T2").Value = 0.03886009588321729
$ws.Range("G3").Value = 23.699655
$ws.Range("H3").Value = 71.09896499999999
$ws.Range("I3").Value = 0.4841969272415696
$ws.Range("J3").Value = 0.4841969272415697
$ws.Range("M3").Value = 7.939250333333333
$ws.Range("O3").Value = 0.7616247559221037
$ws.Range("P3").Value = 0.7616247559221038
$ws.Range("Q3").Value = 188.157493858635
$ws.Range("R3").Value = 1693.417444727715
$ws.Range("S3").Value = 0.368776366528593
$ws.Range("T3").Value = 0.3687763665285931
$ws.Range("G4").Value = 23.699655
$ws.Range("H4").Value = 71.09896499999999
$ws.Range("I4").Value = 0.4841969272415696
$ws.Range("J4").Value = 0.4841969272415697
$ws.Range("M4").Value = 1.648242
$ws.Range("N4").Value = 4.944726
$ws.Range("O4").Value = 0.1581184442163192
$ws.Range("P4").Value = 0.1581184442163192
$ws.Range("Q4").Value = 39.06276675650999
$ws.Range("R4").Value = 351.56490080859
$ws.Range("S4").Value = 0.07656046482975928
$ws.Range("T4").Value = 0.0765604648297593
$ws.Range("H5").Value = 58.032849
$ws.Range("I5").Value = 0.3952142927098025
$ws.Range("J5").Value = 0.3952142927098025
$ws.Range("M5").Value = 0.8366046666666667
$ws.Range("N5").Value = 2.509814
$ws.Range("O5").Value = 0.08025679986157715
$ws.Range("P5").Value = 0.08025679986157715
$ws.Range("Q5").Value = 16.18351743112067
$ws.Range("R5").Value = 145.651656880086
$ws.Range("S5").Value = 0.03171863439244538
$ws.Range("T5").Value = 0.03171863439244538
$ws.Range("H6").Value = 58.032849
$ws.Range("I6").Value = 0.3952142927098025
$ws.Range("J6").Value = 0.3952142927098025
$ws.Range("M6").Value = 7.939250333333333
$ws.Range("O6").Value = 0.7616247559221037
$ws.Range("P6").Value = 0.7616247559221038
$ws.Range("S6").Value = 0.3010049892220301
$ws.Range("T6").Value = 0.3010049892220302
$ws.Range("H7").Value = 58.032849
$ws.Range("I7").Value = 0.3952142927098025
$ws.Range("J7").Value = 0.3952142927098025
$ws.Range("M7").Value = 1.648242
$ws.Range("N7").Value = 4.944726
$ws.Range("O7").Value = 0.1581184442163192
$ws.Range("P7").Value = 0.1581184442163192
$ws.Range("Q7").Value = 31.884059700486
$ws.Range("R7").Value = 286.956537304374
$ws.Range("S7").Value = 0.06249066909532694
$ws.Range("T7").Value = 0.06249066909532695
$ws.Range("G8").Value = 5.902376333333333
$ws.Range("H8").Value = 17.707129
$ws.Range("I8").Value = 0.1205887800486278
$ws.Range("J8").Value = 0.1205887800486278
$ws.Range("M8").Value = 0.8366046666666667
$ws.Range("N8").Value = 2.509814
$ws.Range("O8").Value = 0.08025679986157715
$ws.Range("P8").Value = 0.08025679986157715
$ws.Range("Q8").Value = 4.937955584889555
$ws.Range("R8").Value = 44.441600264006
$ws.Range("S8").Value = 0.009678069585914471
$ws.Range("T8").Value = 0.009678069585914471
$ws.Range("G9").Value = 5.902376333333333
$ws.Range("H9").Value = 17.707129
$ws.Range("I9").Value = 0.1205887800486278
$ws.Range("J9").Value = 0.1205887800486278
$ws.Range("M9").Value = 7.939250333333333
$ws.Range("O9").Value = 0.7616247559221037
$ws.Range("P9").Value = 0.7616247559221038
$ws.Range("Q9").Value = 46.86044327187544
$ws.Range("R9").Value = 421.7439894468789
$ws.Range("S9").Value = 0.09184340017148042
$ws.Range("T9").Value = 0.09184340017148043
$ws.Range("G10").Value = 5.902376333333333
$ws.Range("H10").Value = 17.707129
$ws.Range("I10").Value = 0.1205887800486278
$ws.Range("J10").Value = 0.1205887800486278
$ws.Range("M10").Value = 1.648242
$ws.Range("N10").Value = 4.944726
$ws.Range("O10").Value = 0.1581184442163192
$ws.Range("P10").Value = 0.1581184442163192
$ws.Range("Q10").Value = 9.728544572405999
$ws.Range("R10").Value = 87.55690115165399
$ws.Range("S10").Value = 0.01906731029123294
$ws.Range("T10").Value = 0.01906731029123295
